$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.887.19"
$ws.Range("E2").Value = "  +0.76%  "

$ws.Range("D3").Value = "2.421.03"
$ws.Range("E3").Value = "  +1.01%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "551.09"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").Value = "137.36"
$ws.Range("E6").Value = "  +1.01%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +2.31%  "

$ws.Range("E9").Value = "  -2.09%  "

$ws.Range("D10").Value = "5.67"
$ws.Range("E10").Value = "  -2.77%  "

$ws.Range("E11").Value = "  -2.16%  "

$ws.Range("E12").Value = "  -2.02%  "

$ws.Range("D13").Value = "25.51"
$ws.Range("E13").Value = "  +3.80%  "

$ws.Range("D14").Value = "2.852.99"
$ws.Range("E14").Value = "  +0.73%  "

$ws.Range("D15").Value = "59.822.48"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("E16").Value = "  -1.63%  "

$ws.Range("D17").Value = "2.406.29"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").Value = "11.37"
$ws.Range("E18").Value = "  +1.05%  "

$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D20").Value = "330.82"
$ws.Range("E20").Value = "  -1.48%  "

$ws.Range("D21").Value = "6.71"
$ws.Range("E21").Value = "  -4.30%  "

$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("E23").Value = "  +2.52%  "

$ws.Range("E24").Value = "  +1.30%  "

$ws.Range("E25").Value = "  +3.77%  "

$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("D28").Value = "0.0₃0777"
$ws.Range("E28").Value = "  +2.09%  "

$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("D30").Value = "169.21"
$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("E31").Value = "  -1.77%  "

$ws.Range("D32").Value = "18.68"
$ws.Range("E32").Value = "  -0.21%  "

$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("E35").Value = "  +1.93%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").Value = "4.21"
$ws.Range("E37").Value = "  -1.91%  "

$ws.Range("E38").Value = "  -2.56%  "

$ws.Range("D39").Value = "39.57"
$ws.Range("E39").Value = "  -1.56%  "

$ws.Range("D40").Value = "0.411"
$ws.Range("E40").Value = "  -2.47%  "

$ws.Range("D41").Value = "314.04"
$ws.Range("E41").Value = "  +6.67%  "

$ws.Range("E42").Value = "  -2.02%  "

$ws.Range("D43").Value = "139.31"
$ws.Range("E43").Value = "  -1.68%  "

$ws.Range("D44").Value = "0.0966"
$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("E45").Value = "  -0.50%  "

$ws.Range("D46").Value = "19.57"
$ws.Range("E46").Value = "  +2.68%  "

$ws.Range("E47").Value = "  +1.20%  "

$ws.Range("D48").Value = "0.0224"
$ws.Range("E48").Value = "  -0.56%  "

$ws.Range("D49").Value = "0.392"
$ws.Range("E49").Value = "  -1.56%  "

$ws.Range("E50").Value = "  -0.68%  "

$ws.Range("D51").Value = "11.06"
$ws.Range("E51").Value = "  +0.21%  "
